$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 30 with updated codes/data
$ws.Range("A30").Value = 901
$ws.Range("B30").Value = 381
$ws.Range("C30").Value = 90
$ws.Range("D30").Value = 106
$ws.Range("E30").Value = 97
$ws.Range("F30").Formula = "=+A30-SUM(B30:E30)"

# Update view: scroll back to top and move selection
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("L24").Select()
